{"js": "// Apply the two textual edits described by the diff:\n//  1. In the \"study aims\" paragraph, insert \", specifically, this study will\n//     run over the course of 1 week\" right before \". It would require\n//     participants to be available for 2 in-person meetings in Glasgow\".\n//  2. In the \"financial incentive\" paragraph, drop the stray period right\n//     after \"completed\" and fill in the gift-code amount (\"\u00a3\" -> \"\u00a340\").\n\nconst body = context.document.body;\n\n// --- Edit 1 -------------------------------------------------------------\nconst anchor1 = body.search(\n  \". It would require participants to be available for 2 in-person meetings in Glasgow\",\n  { matchCase: true }\n);\nanchor1.load(\"items\");\nawait context.sync();\n\nif (anchor1.items.length === 0) {\n  throw new Error(\"Could not locate the sentence to insert the new clause before.\");\n}\n\nanchor1.items[0].insertText(\n  \", specifically, this study will run over the course of 1 week\",\n  \"Before\"\n);\nawait context.sync();\n\n// --- Edit 2 -------------------------------------------------------------\nconst anchor2 = body.search(\n  \"completed., with \u00a3 available for full completion.\",\n  { matchCase: true }\n);\nanchor2.load(\"items\");\nawait context.sync();\n\nif (anchor2.items.length === 0) {\n  throw new Error(\"Could not locate the incentive sentence to update.\");\n}\n\nanchor2.items[0].insertText(\n  \"completed, with \u00a340 available for full completion.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Apply the two textual edits described by the diff:\n#  1. In the \"study aims\" paragraph, insert \", specifically, this study will\n#     run over the course of 1 week\" right before \". It would require\n#     participants to be available for 2 in-person meetings in Glasgow\".\n#  2. In the \"financial incentive\" paragraph, drop the stray period right\n#     after \"completed\" and fill in the gift-code amount (\"\u00a3\" -> \"\u00a340\").\n\n$doc = $word.ActiveDocument\n\n# --- Edit 1 ---------------------------------------------------------------\n$r1 = $doc.Content\n$found1 = $r1.Find.Execute(\". It would require participants to be available for 2 in-person meetings in Glasgow\")\nif (-not $found1) {\n    throw \"Could not locate the sentence to insert the new clause before.\"\n}\n$r1.Collapse(1)  # wdCollapseStart\n$r1.InsertBefore(\", specifically, this study will run over the course of 1 week\")\n\n# --- Edit 2 ---------------------------------------------------------------\n$r2 = $doc.Content\n$found2 = $r2.Find.Execute(\"completed., with \u00a3 available for full completion.\")\nif (-not $found2) {\n    throw \"Could not locate the incentive sentence to update.\"\n}\n$r2.Text = \"completed, with \u00a340 available for full completion.\"\n"}
